# Applies the "Updated cryptos list" data refresh to the cryptocurrency table on the active sheet.
# For each changed cell we set the new text value; for Price (column D) cells whose new value
# looks like a plain number (e.g. "1.00", "0.0214"), we first force the cells number format to
# Text ("@") so Excel keeps it as a literal string instead of silently converting it to a number
# (which would drop meaningful trailing zeros / change formatting, e.g. "1.00" -> 1, "6.10" -> 6.1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '56.646.91'
$ws.Range('E2').Value = '  +0.02%  '
# Row 3
$ws.Range('D3').Value = '2.320.54'
$ws.Range('E3').Value = '  -0.11%  '
# Row 4
$ws.Range('E4').Value = '  -0.05%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '514.81'
$ws.Range('E5').Value = '  -1.05%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.89'
$ws.Range('E6').Value = '  -2.30%  '
# Row 7
$ws.Range('E7').Value = '  +0.33%  '
# Row 8
$ws.Range('E8').Value = '  -0.54%  '
# Row 9
$ws.Range('E9').Value = '  -2.75%  '
# Row 10
$ws.Range('E10').Value = '  -0.07%  '
# Row 11
$ws.Range('E11').Value = '  -0.89%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.337'
$ws.Range('E12').Value = '  -1.63%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.62'
$ws.Range('E13').Value = '  -1.53%  '
# Row 14
$ws.Range('D14').Value = '2.734.59'
$ws.Range('E14').Value = '  -0.46%  '
# Row 15
$ws.Range('D15').Value = '56.620.43'
$ws.Range('E15').Value = '  -0.19%  '
# Row 16
$ws.Range('E16').Value = '  -1.02%  '
# Row 17
$ws.Range('D17').Value = '2.329.62'
$ws.Range('E17').Value = '  -0.24%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.37'
$ws.Range('E18').Value = '  -1.28%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '328.24'
$ws.Range('E19').Value = '  +1.75%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.16'
$ws.Range('E20').Value = '  -1.69%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.74'
$ws.Range('E21').Value = '  +1.99%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.14%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.24'
$ws.Range('E23').Value = '  +1.10%  '
# Row 24
$ws.Range('E24').Value = '  -1.03%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.62'
$ws.Range('E25').Value = '  +8.19%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  +0.09%  '
# Row 27
$ws.Range('E27').Value = '  +1.61%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '167.42'
$ws.Range('E28').Value = '  +0.08%  '
# Row 29
$ws.Range('E29').Value = '  -2.43%  '
# Row 30
$ws.Range('D30').Value = '0.0₃0720'
$ws.Range('E30').Value = '  -3.08%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.10'
$ws.Range('E31').Value = '  -1.85%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.29'
$ws.Range('E32').Value = '  -0.31%  '
# Row 33
$ws.Range('E33').Value = '  +0.00%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.997'
$ws.Range('E34').Value = '  +0.48%  '
# Row 35
$ws.Range('E35').Value = '  -1.19%  '
# Row 36
$ws.Range('E36').Value = '  -2.54%  '
# Row 37
$ws.Range('E37').Value = '  -4.59%  '
# Row 38
$ws.Range('E38').Value = '  +0.57%  '
# Row 39
$ws.Range('E39').Value = '  +1.80%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '149.30'
$ws.Range('E40').Value = '  +7.24%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.374'
$ws.Range('E41').Value = '  -1.44%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.57'
$ws.Range('E42').Value = '  -0.96%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '276.44'
$ws.Range('E43').Value = '  -0.34%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.01'
$ws.Range('E44').Value = '  -4.91%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0928'
$ws.Range('E45').Value = '  -0.56%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0495'
$ws.Range('E46').Value = '  -2.19%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.557'
$ws.Range('E47').Value = '  -1.30%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.26'
$ws.Range('E48').Value = '  +2.21%  '
# Row 49
$ws.Range('B49').Value = 'Polygon'
$ws.Range('C49').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.380'
$ws.Range('E49').Value = '  +0.15%  '
# Row 50
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0214'
$ws.Range('E50').Value = '  -1.57%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.09'
$ws.Range('E51').Value = '  +1.25%  '
